$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a brand-new paragraph right after "...organize and arrange your
#    work. " (paragraph 3) and before "The inspiration from the W3Schools..."
#    (paragraph 4). This new paragraph carries the "Tree View is unique..."
#    sentence that used to be appended onto the inspiration paragraph.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "Tree View is unique in that it does not implement other code snippets from W3Schools, but instead integrates with your data. Or other sources of data given the WebAPI feature."

# ---------------------------------------------------------------------------
# 2) Strip the now-duplicated tail off the old "inspiration" paragraph,
#    leaving just the period+space right after the hyperlink. We deliberately
#    leave the leading ". " untouched so its run keeps its original
#    (non-hyperlink) character formatting.
# ---------------------------------------------------------------------------
$old = "It is unique in that it does not implement another code snippet from W3Schools, but instead integrates with your data. Or other sources of data given the WebAPI feature. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Relocate the "_GoBack" bookmark: it used to sit right before "feature."
#    near the end of the old combined paragraph; now it belongs right before
#    "other code snippet" in the new paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$found = $d.Content
$found.Find.Execute("other code snippet", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPos = $d.Range($found.Start, $found.Start)
$d.Bookmarks.Add("_GoBack", $bmPos)

Write-Host "Paragraphs count: " $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para $i : [$($p.Range.Text)]"
}
